# Update "paises.xlsx" worksheet with the latest COVID-19 country figures
# (data refresh + India/Irlanda/Chile and Kazajistan/Crucero re-sort by
# "Casos totales").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (row 1) ---------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 8 de Abril de 2020 a las 17:52"

# --- Estados Unidos (row 4) -------------------------------------------------
$ws.Range("B4").Value = 404156
$ws.Range("C4").Value = 3821
$ws.Range("E4").Value = 369353
$ws.Range("F4").Value = 9220

# --- Iran (row 10) -----------------------------------------------------------
$ws.Range("D10").Value = 29812
$ws.Range("E10").Value = 30781

# --- Canada (row 16) ---------------------------------------------------------
$ws.Range("B16").Value = 18479
$ws.Range("C16").Value = 582
$ws.Range("D16").Value = 4333
$ws.Range("E16").Value = 13744
$ws.Range("G16").Value = 21
$ws.Range("H16").Value = 402

# --- India overtakes Irlanda and Chile (rows 26-28) --------------------------
# Row 26 becomes India with freshly updated figures.
$ws.Range("A26").Value = "India"
$ws.Range("B26").Value = 5749
$ws.Range("C26").Value = 398
$ws.Range("D26").Value = 468
$ws.Range("E26").Value = 5117
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = 164

# Row 27 becomes Irlanda (its previous figures, now one row lower).
$ws.Range("A27").Value = "Irlanda"
$ws.Range("B27").Value = 5709
$ws.Range("C27").Value = 0
$ws.Range("D27").Value = 25
$ws.Range("E27").Value = 5474
$ws.Range("F27").Value = 165
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 210

# Row 28 becomes Chile (its previous figures, now one row lower).
$ws.Range("A28").Value = "Chile"
$ws.Range("B28").Value = 5546
$ws.Range("C28").Value = 430
$ws.Range("D28").Value = 1115
$ws.Range("E28").Value = 4383
$ws.Range("F28").Value = 362
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = 48

# --- Pakistan (row 34) --------------------------------------------------------
$ws.Range("B34").Value = 4196
$ws.Range("C34").Value = 161
$ws.Range("E34").Value = 3669

# --- Luxemburgo (row 38) ------------------------------------------------------
$ws.Range("B38").Value = 3034
$ws.Range("C38").Value = 64
$ws.Range("E38").Value = 2488
$ws.Range("F38").Value = 34
$ws.Range("G38").Value = 2
$ws.Range("H38").Value = 46

# --- Kazajistan overtakes Crucero (rows 75-76) --------------------------------
# Row 75 becomes Kazajistan with freshly updated figures.
$ws.Range("A75").Value = "Kazajistan"
$ws.Range("B75").Value = 718
$ws.Range("C75").Value = 21
$ws.Range("D75").Value = 54
$ws.Range("E75").Value = 657
$ws.Range("F75").Value = 21
$ws.Range("G75").Value = 1
$ws.Range("H75").Value = 7

# Row 76 becomes Crucero (its previous figures, now one row lower).
$ws.Range("A76").Value = "Crucero"
$ws.Range("B76").Value = 712
$ws.Range("C76").Value = 0
$ws.Range("D76").Value = 619
$ws.Range("E76").Value = 82
$ws.Range("F76").Value = 10
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 11

# --- Tanzania (row 160) -------------------------------------------------------
$ws.Range("B160").Value = 25
$ws.Range("C160").Value = 1
$ws.Range("E160").Value = 19
